# Apply cryptos list update (prices / 1h volume %) per commit
# "Updated cryptos list on Tue Aug 20 14:36:52 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit per-cell updates. For numeric-looking text values in the
# "Price" column (D), a leading apostrophe forces Excel to keep them as text
# (preserving formatting such as trailing zeros / grouping dots) instead of
# auto-converting to a number; the style is then reset to Normal so no new
# number-format styling is left on the cell.

$ws.Range("D2").Value = "59.682.80"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "2.597.90"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'569.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.22%  "
$ws.Range("D6").Value = "'143.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").Value = "2.604.46"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'6.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("E12").Value = "  +10.42%  "
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("D14").Value = "3.062.07"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "59.678.74"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").Value = "'21.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.18%  "
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "2.615.59"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "'4.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "'336.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "'10.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").Value = "'6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'65.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "'0.447"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.07%  "
$ws.Range("D28").Value = "'7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").Value = "  +7.75%  "
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "'159.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "'18.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'4.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("E36").Value = "  +10.00%  "
$ws.Range("D37").Value = "'1.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").Value = "'0.869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").Value = "'37.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("D41").Value = "'293.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'0.0976"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "'0.593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "'10.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "'19.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("D49").Value = "'125.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.26%  "
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("D51").Value = "1.936.83"
$ws.Range("E51").Value = "  +1.78%  "
